# Applies the cryptos-list refresh described in the commit:
# "Updated cryptos list on Sun Nov 17 05:32:10 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a Price (column D) value while forcing it to stay plain text
# (many prices look numeric, e.g. "1.00" / "45.00" / "14.30", and a bare
# Range.Value assignment would have Excel coerce them into numbers and drop
# the trailing zeros -- the workbook stores these as inline strings).
# Resetting the Style back to "Normal" afterwards keeps the cell from
# picking up a lingering explicit number-format style.
function Set-PriceText($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
Set-PriceText $ws.Range("D2") "90.684.33"
$ws.Range("E2").Value = "  -0.97%  "

# Row 3
Set-PriceText $ws.Range("D3") "3.111.89"
$ws.Range("E3").Value = "  -0.61%  "

# Row 4
Set-PriceText $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.16%  "

# Row 5
Set-PriceText $ws.Range("D5") "236.82"
$ws.Range("E5").Value = "  +7.24%  "

# Row 6
Set-PriceText $ws.Range("D6") "626.49"
$ws.Range("E6").Value = "  +0.49%  "

# Row 7
Set-PriceText $ws.Range("D7") "1.03"
$ws.Range("E7").Value = "  +6.94%  "

# Row 8
Set-PriceText $ws.Range("D8") "0.354"
$ws.Range("E8").Value = "  -5.98%  "

# Row 9
$ws.Range("E9").Value = "  +0.06%  "

# Row 10
Set-PriceText $ws.Range("D10") "3.381.81"
$ws.Range("E10").Value = "  +8.12%  "

# Row 11
Set-PriceText $ws.Range("D11") "0.714"
$ws.Range("E11").Value = "  -4.10%  "

# Row 12
$ws.Range("E12").Value = "  +3.26%  "

# Row 13
Set-PriceText $ws.Range("D13") "36.35"
$ws.Range("E13").Value = "  +4.02%  "

# Row 14
$ws.Range("E14").Value = "  -3.60%  "

# Row 16
Set-PriceText $ws.Range("D16") "90.325.55"
$ws.Range("E16").Value = "  -0.89%  "

# Row 17
Set-PriceText $ws.Range("D17") "3.685.01"
$ws.Range("E17").Value = "  -0.62%  "

# Row 18
Set-PriceText $ws.Range("D18") "3.065.73"
$ws.Range("E18").Value = "  -1.72%  "

# Row 19
Set-PriceText $ws.Range("D19") "3.67"
$ws.Range("E19").Value = "  -2.66%  "

# Row 20
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-PriceText $ws.Range("D20") "14.30"
$ws.Range("E20").Value = "  +1.04%  "

# Row 21
$ws.Range("B21").Value = "PEPE"
$ws.Range("C21").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-PriceText $ws.Range("D21") "0.0000213"
$ws.Range("E21").Value = "  -6.02%  "

# Row 22
Set-PriceText $ws.Range("D22") "449.30"
$ws.Range("E22").Value = "  +3.22%  "

# Row 23
Set-PriceText $ws.Range("D23") "8.96"
$ws.Range("E23").Value = "  +1.53%  "

# Row 24
$ws.Range("E24").Value = "  +3.80%  "

# Row 25
Set-PriceText $ws.Range("D25") "6.15"
$ws.Range("E25").Value = "  +0.97%  "

# Row 26
Set-PriceText $ws.Range("D26") "90.40"
$ws.Range("E26").Value = "  +5.11%  "

# Row 27
Set-PriceText $ws.Range("D27") "12.25"
$ws.Range("E27").Value = "  -1.58%  "

# Row 28
$ws.Range("E28").Value = "  -0.63%  "

# Row 29
Set-PriceText $ws.Range("D29") "1.00"
$ws.Range("E29").Value = "  +0.05%  "

# Row 30
$ws.Range("E30").Value = "  +3.46%  "

# Row 31
Set-PriceText $ws.Range("D31") "0.159"
$ws.Range("E31").Value = "  -6.67%  "

# Row 32
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-PriceText $ws.Range("D32") "27.47"
$ws.Range("E32").Value = "  +15.78%  "

# Row 33
$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-PriceText $ws.Range("D33") "0.195"
$ws.Range("E33").Value = "  +27.12%  "

# Row 34
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-PriceText $ws.Range("D34") "0.151"
$ws.Range("E34").Value = "  +3.89%  "

# Row 35
$ws.Range("B35").Value = "dogwifhat"
$ws.Range("C35").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-PriceText $ws.Range("D35") "3.75"
$ws.Range("E35").Value = "  -2.73%  "

# Row 36
Set-PriceText $ws.Range("D36") "507.98"
$ws.Range("E36").Value = "  -4.76%  "

# Row 37
$ws.Range("B37").Value = "PancakeSwap"
$ws.Range("C37").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-PriceText $ws.Range("D37") "1.92"
$ws.Range("E37").Value = "  +3.11%  "

# Row 38
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-PriceText $ws.Range("D38") "7.00"
$ws.Range("E38").Value = "  -3.41%  "

# Row 39
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-PriceText $ws.Range("D39") "1.31"
$ws.Range("E39").Value = "  +0.60%  "

# Row 40
$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-PriceText $ws.Range("D40") "0.423"
$ws.Range("E40").Value = "  +11.01%  "

# Row 41
$ws.Range("B41").Value = "WhiteBITCoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-PriceText $ws.Range("D41") "22.19"
$ws.Range("E41").Value = "  -0.60%  "

# Row 42
$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-PriceText $ws.Range("D42") "0.0854"
$ws.Range("E42").Value = "  +9.11%  "

# Row 43
$ws.Range("E43").Value = "  +0.02%  "

# Row 44
$ws.Range("B44").Value = "Binance-PegBSC-USD"
$ws.Range("C44").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-PriceText $ws.Range("D44") "0.738"
$ws.Range("E44").Value = "  -17.69%  "

# Row 45
Set-PriceText $ws.Range("D45") "3.25"
$ws.Range("E45").Value = "  +33.01%  "

# Row 46
$ws.Range("E46").Value = "  +1.11%  "

# Row 47
Set-PriceText $ws.Range("D47") "0.697"
$ws.Range("E47").Value = "  +11.47%  "

# Row 48
Set-PriceText $ws.Range("D48") "149.24"
$ws.Range("E48").Value = "  +2.68%  "

# Row 49
Set-PriceText $ws.Range("D49") "4.54"
$ws.Range("E49").Value = "  +8.15%  "

# Row 50
$ws.Range("B50").Value = "OKB"
$ws.Range("C50").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-PriceText $ws.Range("D50") "45.00"
$ws.Range("E50").Value = "  +1.88%  "

# Row 51
$ws.Range("B51").Value = "ImmutableX"
$ws.Range("C51").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-PriceText $ws.Range("D51") "1.35"
$ws.Range("E51").Value = "  +3.97%  "
